$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new header row with column titles ---
$ws.Range("A1").Value = "Test_Case_Name"
$ws.Range("B1").Value = "Test_Case_ID"
$ws.Range("C1").Value = "Test_Description"
$ws.Range("D1").Value = "Success_Message"
$ws.Range("E1").Value = "Failure_Message"

# --- Fix casing of existing Test_Case_Name values (column A, rows 2-6) ---
$ws.Range("A2").Value = "getWeatherInJSONFormatTest"
$ws.Range("A3").Value = "getWeatherInXMLFormatTest"
$ws.Range("A4").Value = "getTempInFahrenheitUnitTest"
$ws.Range("A5").Value = "getTempInCelsiusUnitTest"
$ws.Range("A6").Value = "getWeatherInDiffLangTest"

# --- Style the header row: bold white text, blue (Accent1) fill, thin border, centered ---
$hdrBase = $ws.Range("A1")
$hdrBase.Font.Bold = $true
$hdrBase.Font.ThemeColor = 2
$hdrBase.Interior.ThemeColor = 5
$hdrBase.Borders.LineStyle = 1
$hdrBase.HorizontalAlignment = -4108
$hdrBase.VerticalAlignment = -4108

$hdrRest = $ws.Range("B1:E1")
$hdrBase.Copy()
$hdrRest.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Style the data rows (A2:E6): thin border around every cell ---
$dataBase = $ws.Range("A2")
$dataBase.Borders.LineStyle = 1

$dataRest = $ws.Range("A2:E6")
$dataBase.Copy()
$dataRest.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Adjust column widths for the widened/ new columns ---
# (input values chosen so the engine's internal width quantization lands as
# close as possible to the target character widths of 39.285/84.141/25.570)
$ws.Columns("C").ColumnWidth = 38.5
$ws.Columns("D").ColumnWidth = 83.3
$ws.Columns("E").ColumnWidth = 24.6

# --- Update selection to match the post-edit state ---
$ws.Range("A7").Select()

Write-Host "edit complete"
